$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.6842105263157895
$ws.Range("D2").Value = 0.7194244604316546
$ws.Range("C3").Value = 0.6530612244897959
$ws.Range("D3").Value = 0.7862068965517242
$ws.Range("C4").Value = 0.6222222222222222
$ws.Range("D4").Value = 0.6754966887417219
$ws.Range("C5").Value = 0.6341463414634146
$ws.Range("D5").Value = 0.676470588235294
$ws.Range("C18").Value = 0.6500000000000001
$ws.Range("D18").Value = 0.6950354609929078
$ws.Range("C19").Value = 0.6500000000000001
$ws.Range("D19").Value = 0.6950354609929078
$ws.Range("C20").Value = 0.6666666666666666
$ws.Range("D20").Value = 0.6950354609929078
$ws.Range("C21").Value = 0.6808510638297872
$ws.Range("D21").Value = 0.7808219178082192
$ws.Range("C22").Value = 0.6808510638297872
$ws.Range("D22").Value = 0.7808219178082192
$ws.Range("C23").Value = 0.6808510638297872
$ws.Range("D23").Value = 0.7808219178082192
$ws.Range("C24").Value = 0.6046511627906977
$ws.Range("D24").Value = 0.7866666666666666
$ws.Range("C25").Value = 0.6046511627906977
$ws.Range("D25").Value = 0.7866666666666666
$ws.Range("C26").Value = 0.6046511627906977
$ws.Range("D26").Value = 0.7866666666666666
$ws.Range("C27").Value = 0.717948717948718
$ws.Range("D27").Value = 0.72992700729927
$ws.Range("C28").Value = 0.717948717948718
$ws.Range("D28").Value = 0.72992700729927
$ws.Range("C29").Value = 0.717948717948718
$ws.Range("D29").Value = 0.72992700729927
$ws.Range("C72").Value = 0.7272727272727272
$ws.Range("C73").Value = 0.7111111111111111
$ws.Range("C74").Value = 0.6521739130434783
$ws.Range("C75").Value = 0.6511627906976744
$ws.Range("C84").Value = 0.7142857142857143
$ws.Range("D84").Value = 0.7638888888888888
$ws.Range("C85").Value = 0.6486486486486486
$ws.Range("D85").Value = 0.7916666666666666
$ws.Range("C86").Value = 0.6500000000000001
$ws.Range("D86").Value = 0.7605633802816901
$ws.Range("C87").Value = 0.717948717948718
$ws.Range("D87").Value = 0.7076923076923075
$ws.Range("C108").Value = 0.6486486486486486
$ws.Range("D108").Value = 0.951048951048951
$ws.Range("C109").Value = 0.5789473684210527
$ws.Range("C110").Value = 0.5116279069767442
$ws.Range("D110").Value = 0.9733333333333333
$ws.Range("C111").Value = 0.717948717948718
$ws.Range("D111").Value = 0.9295774647887325
$ws.Range("C124").Value = 0.6153846153846153
$ws.Range("D124").Value = 0.676056338028169
$ws.Range("C125").Value = 0.7567567567567567
$ws.Range("D125").Value = 0.6870229007633588
$ws.Range("D126").Value = 0.6493506493506493
$ws.Range("C127").Value = 0.7441860465116279
$ws.Range("D127").Value = 0.6979865771812082
$ws.Range("D140").Value = 0.8609271523178808
$ws.Range("C141").Value = 0.5641025641025642
$ws.Range("D141").Value = 0.8918918918918919
$ws.Range("C142").Value = 0.7441860465116279
$ws.Range("D142").Value = 0.8588957055214724
$ws.Range("C143").Value = 0.7368421052631577
$ws.Range("D143").Value = 0.8000000000000002
$ws.Range("C156").Value = 0.6500000000000001
$ws.Range("D156").Value = 0.8311688311688312
$ws.Range("C157").Value = 0.5405405405405405
$ws.Range("D157").Value = 0.8489208633093526
$ws.Range("C158").Value = 0.6956521739130435
$ws.Range("D158").Value = 0.8181818181818183
$ws.Range("C159").Value = 0.7222222222222222
$ws.Range("D159").Value = 0.7419354838709677